$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the site-name header/value from column A to new column J
$ws.Range("J1").Value = $ws.Range("A1").Value2
$ws.Range("J2").Value = $ws.Range("A2").Value2

# Clear old column A entirely, then set A2 to numeric 0
$ws.Range("A1:A2").Clear()
$ws.Range("A2").Value = 0

$ws.Range("D6").Select()
